$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.506.96"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.843.12"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D5").Value = "'262.09"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.5319"
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("D8").Value = "'0.3069"
$ws.Range("E8").Value = "  -5.01%  "
$ws.Range("D9").Value = "'0.06901"
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("D10").Value = "'18.37"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").Value = "'0.07804"
$ws.Range("D12").Value = "'0.7537"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").Value = "1.839.98"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "'89.70"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "'5.011"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "'14.00"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "'0.000007948"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "26.524.83"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "'4.618"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "'5.982"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").Value = "'9.316"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "'142.53"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "'2.192"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").Value = "'1.688"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").Value = "'17.01"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").Value = "'111.29"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").Value = "'4.270"
$ws.Range("E29").Value = "  +2.62%  "
$ws.Range("D30").Value = "'0.08804"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").Value = "'4.084"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").Value = "'0.04815"
$ws.Range("D33").Value = "'2.930"
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("D34").Value = "'0.7290"
$ws.Range("E34").Value = "  +2.06%  "
$ws.Range("D35").Value = "'1.133"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").Value = "'2.306"
$ws.Range("E37").Value = "  +5.81%  "
$ws.Range("D38").Value = "'0.01716"
$ws.Range("E38").Value = "  -3.80%  "
$ws.Range("D39").Value = "'0.4795"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").Value = "'0.9068"
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("D41").Value = "'108.18"
$ws.Range("E41").Value = "  -3.73%  "
$ws.Range("D42").Value = "'5.866"
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "'7.499"
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.4136"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'9.032"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").Value = "'0.1241"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").Value = "'0.9004"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'34.85"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").Value = "'60.24"
$ws.Range("E51").Value = "  +0.64%  "

